# Apply the row-rotation edits for rows 15-18 and 25-27 on sheet "Artfynd".
# The underlying data rows were re-ordered (each row's content shifted to
# an adjacent row number); we reproduce the resulting per-cell values directly.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artfynd")

# Row 15
$ws.Cells.Item(15, 1).Value = 111732650
$ws.Cells.Item(15, 2).Value = 88489
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "NT"
$ws.Cells.Item(15, 5).Value = 1962
$ws.Cells.Item(15, 6).NumberFormat = "@"
$ws.Cells.Item(15, 6).Value = "Vaddporing"
$ws.Cells.Item(15, 7).NumberFormat = "@"
$ws.Cells.Item(15, 7).Value = "Anomoporia kamtschatica"
$ws.Cells.Item(15, 8).NumberFormat = "@"
$ws.Cells.Item(15, 8).Value = "(Parmasto) Bondartseva"
$ws.Cells.Item(15, 10).Value = ""
$ws.Cells.Item(15, 17).Value = 367499.9406743076
$ws.Cells.Item(15, 18).Value = 6871172.113255707
$ws.Cells.Item(15, 25).NumberFormat = "@"
$ws.Cells.Item(15, 25).Value = "2023-07-29"
$ws.Cells.Item(15, 26).NumberFormat = "@"
$ws.Cells.Item(15, 26).Value = "12:00"
$ws.Cells.Item(15, 27).NumberFormat = "@"
$ws.Cells.Item(15, 27).Value = "2023-07-29"
$ws.Cells.Item(15, 28).NumberFormat = "@"
$ws.Cells.Item(15, 28).Value = "12:00"
$ws.Cells.Item(15, 29).NumberFormat = "@"
$ws.Cells.Item(15, 29).Value = "Växer under silverved"
$ws.Cells.Item(15, 35).NumberFormat = "@"
$ws.Cells.Item(15, 35).Value = "Sandtallskog. Kontinuitetskog"

# Row 16
$ws.Cells.Item(16, 1).Value = 111729509
$ws.Cells.Item(16, 10).NumberFormat = "@"
$ws.Cells.Item(16, 10).Value = "fruktkroppar"
$ws.Cells.Item(16, 17).Value = 367351.5527010285
$ws.Cells.Item(16, 18).Value = 6871697.886848727
$ws.Cells.Item(16, 26).NumberFormat = "@"
$ws.Cells.Item(16, 26).Value = "14:30"
$ws.Cells.Item(16, 28).NumberFormat = "@"
$ws.Cells.Item(16, 28).Value = "14:30"

# Row 17
$ws.Cells.Item(17, 1).Value = 111730480
$ws.Cells.Item(17, 2).Value = 90666
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "LC"
$ws.Cells.Item(17, 5).Value = 4364
$ws.Cells.Item(17, 6).NumberFormat = "@"
$ws.Cells.Item(17, 6).Value = "Dropptaggsvamp"
$ws.Cells.Item(17, 7).NumberFormat = "@"
$ws.Cells.Item(17, 7).Value = "Hydnellum ferrugineum"
$ws.Cells.Item(17, 8).NumberFormat = "@"
$ws.Cells.Item(17, 8).Value = "(Fr.:Fr.) P. Karst."
$ws.Cells.Item(17, 10).Value = ""
$ws.Cells.Item(17, 17).Value = 367528.9347990834
$ws.Cells.Item(17, 18).Value = 6871323.127923099

# Row 18
$ws.Cells.Item(18, 1).Value = 111730457
$ws.Cells.Item(18, 2).Value = 90660
$ws.Cells.Item(18, 5).Value = 4362
$ws.Cells.Item(18, 6).NumberFormat = "@"
$ws.Cells.Item(18, 6).Value = "Blå taggsvamp"
$ws.Cells.Item(18, 7).NumberFormat = "@"
$ws.Cells.Item(18, 7).Value = "Hydnellum caeruleum"
$ws.Cells.Item(18, 8).NumberFormat = "@"
$ws.Cells.Item(18, 8).Value = "(Hornem.) P.Karst."
$ws.Cells.Item(18, 10).NumberFormat = "@"
$ws.Cells.Item(18, 10).Value = "fruktkroppar"
$ws.Cells.Item(18, 17).Value = 367515.5716414675
$ws.Cells.Item(18, 18).Value = 6871294.915694831
$ws.Cells.Item(18, 25).NumberFormat = "@"
$ws.Cells.Item(18, 25).Value = "2023-08-27"
$ws.Cells.Item(18, 26).NumberFormat = "@"
$ws.Cells.Item(18, 26).Value = "13:30"
$ws.Cells.Item(18, 27).NumberFormat = "@"
$ws.Cells.Item(18, 27).Value = "2023-08-27"
$ws.Cells.Item(18, 28).NumberFormat = "@"
$ws.Cells.Item(18, 28).Value = "13:30"
$ws.Cells.Item(18, 29).Value = ""
$ws.Cells.Item(18, 35).NumberFormat = "@"
$ws.Cells.Item(18, 35).Value = "Sandtallskog. Kontinuitetsskog"

# Row 25
$ws.Cells.Item(25, 1).Value = 111730155
$ws.Cells.Item(25, 2).Value = 90682
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "NT"
$ws.Cells.Item(25, 5).Value = 2059
$ws.Cells.Item(25, 6).NumberFormat = "@"
$ws.Cells.Item(25, 6).Value = "Skrovlig taggsvamp"
$ws.Cells.Item(25, 7).NumberFormat = "@"
$ws.Cells.Item(25, 7).Value = "Hydnellum scabrosum"
$ws.Cells.Item(25, 8).NumberFormat = "@"
$ws.Cells.Item(25, 8).Value = "(Fr.) E.Larss., K.H.Larss. & Kõljalg"
$ws.Cells.Item(25, 9).NumberFormat = "@"
$ws.Cells.Item(25, 9).Value = "15"
$ws.Cells.Item(25, 17).Value = 367511.4922534205
$ws.Cells.Item(25, 18).Value = 6871287.067064899
$ws.Cells.Item(25, 25).NumberFormat = "@"
$ws.Cells.Item(25, 25).Value = "2023-08-27"
$ws.Cells.Item(25, 26).NumberFormat = "@"
$ws.Cells.Item(25, 26).Value = "00:00"
$ws.Cells.Item(25, 27).NumberFormat = "@"
$ws.Cells.Item(25, 27).Value = "2023-08-27"
$ws.Cells.Item(25, 28).NumberFormat = "@"
$ws.Cells.Item(25, 28).Value = "00:00"

# Row 26
$ws.Cells.Item(26, 1).Value = 111733035
$ws.Cells.Item(26, 2).Value = 8377
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "LC"
$ws.Cells.Item(26, 5).Value = 106545
$ws.Cells.Item(26, 6).NumberFormat = "@"
$ws.Cells.Item(26, 6).Value = "Mindre märgborre"
$ws.Cells.Item(26, 7).NumberFormat = "@"
$ws.Cells.Item(26, 7).Value = "Tomicus minor"
$ws.Cells.Item(26, 8).NumberFormat = "@"
$ws.Cells.Item(26, 8).Value = "(Hartig, 1834)"
$ws.Cells.Item(26, 9).Value = ""
$ws.Cells.Item(26, 10).Value = ""
$ws.Cells.Item(26, 12).Value = ""
$ws.Cells.Item(26, 13).NumberFormat = "@"
$ws.Cells.Item(26, 13).Value = "äldre gnagspår"
$ws.Cells.Item(26, 17).Value = 367525.4248074447
$ws.Cells.Item(26, 18).Value = 6871378.373665834
$ws.Cells.Item(26, 25).NumberFormat = "@"
$ws.Cells.Item(26, 25).Value = "2023-07-29"
$ws.Cells.Item(26, 26).NumberFormat = "@"
$ws.Cells.Item(26, 26).Value = "11:30"
$ws.Cells.Item(26, 27).NumberFormat = "@"
$ws.Cells.Item(26, 27).Value = "2023-07-29"
$ws.Cells.Item(26, 28).NumberFormat = "@"
$ws.Cells.Item(26, 28).Value = "11:30"

# Row 27
$ws.Cells.Item(27, 1).Value = 111730937
$ws.Cells.Item(27, 2).Value = 88032
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "VU"
$ws.Cells.Item(27, 5).Value = 6276
$ws.Cells.Item(27, 6).NumberFormat = "@"
$ws.Cells.Item(27, 6).Value = "Goliatmusseron"
$ws.Cells.Item(27, 7).NumberFormat = "@"
$ws.Cells.Item(27, 7).Value = "Tricholoma matsutake"
$ws.Cells.Item(27, 8).NumberFormat = "@"
$ws.Cells.Item(27, 8).Value = "(S.Ito & S.Imai) Singer"
$ws.Cells.Item(27, 9).NumberFormat = "@"
$ws.Cells.Item(27, 9).Value = "1"
$ws.Cells.Item(27, 10).NumberFormat = "@"
$ws.Cells.Item(27, 10).Value = "fruktkroppar"
$ws.Cells.Item(27, 12).Value = ""
$ws.Cells.Item(27, 13).Value = ""
$ws.Cells.Item(27, 17).Value = 367427.0662824844
$ws.Cells.Item(27, 18).Value = 6871565.895031672
$ws.Cells.Item(27, 25).NumberFormat = "@"
$ws.Cells.Item(27, 25).Value = "2023-08-25"
$ws.Cells.Item(27, 26).NumberFormat = "@"
$ws.Cells.Item(27, 26).Value = "13:00"
$ws.Cells.Item(27, 27).NumberFormat = "@"
$ws.Cells.Item(27, 27).Value = "2023-08-25"
$ws.Cells.Item(27, 28).NumberFormat = "@"
$ws.Cells.Item(27, 28).Value = "13:00"
